$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new formula in cell C11 (row 11)
$ws.Range("C11").Formula = "=12.3 - 7.2"

# Update the selection to match the post-edit state (C12)
$ws.Range("C12").Select()
